$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra (blank) row 2 that separated the header from the data,
# shifting everything below it up by one row.
$ws.Rows("2:2").Delete()

# Restore the active cell selection to match the post-edit layout.
$ws.Range("B5").Select() | Out-Null
